$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# Updated "GDP per Capita" values for Djibouti, years 1950-2016 (row 2 = 1950 ... row 68 = 2016).
# This supersedes the previous 1950-2008 series with refreshed figures and extends it through 2016.
$gdpPerCapita = @(
    "3459",
    "3566",
    "3585",
    "3633",
    "3743",
    "3765",
    "3810",
    "3848",
    "3840",
    "3948",
    "4087",
    "4114",
    "4057",
    "4092",
    "4055",
    "4047",
    "4063",
    "4055",
    "4028",
    "4018",
    "4774",
    "4846",
    "4787",
    "4801",
    "4414",
    "4691",
    "4878",
    "4036",
    "3419",
    "3308",
    "3167",
    "3091",
    "3043",
    "2963",
    "2944",
    "2880",
    "2743",
    "2670",
    "2652",
    "2436",
    "2308",
    "2013.34268624068",
    "2141.63413537219",
    "2290.78463354907",
    "2175.36187332938",
    "2262.42871519311",
    "2081.99546483789",
    "2080.51341346604",
    "1974.56567647215",
    "1930.50779212697",
    "1880.99875054224",
    "1845.69908201712",
    "1824.45508650037",
    "1945.66548268343",
    "2097.12337871016",
    "2126.80514473055",
    "2342.21916058477",
    "2396.01808306678",
    "2463.44491789316",
    "2521.50471546992",
    "2543.87329499783",
    "2590",
    "2654",
    "2724",
    "2823",
    "2941",
    "3064"
)

$firstRow = 2
$lastRow = $firstRow + $gdpPerCapita.Length - 1

# Keep the Data column stored as text (matching the workbook's existing convention).
$ws.Range("E$firstRow`:E$lastRow").NumberFormat = "@"

for ($i = 0; $i -lt $gdpPerCapita.Length; $i++) {
    $row = $firstRow + $i
    $year = 1950 + $i

    # Rows 61-68 (years 2009-2016) are brand new; populate their identifying columns too.
    if ($row -gt 60) {
        $ws.Cells.Item($row, 1).Value = 262.0
        $ws.Cells.Item($row, 2).Value = "Djibouti"
        $ws.Cells.Item($row, 3).Value = "GDP per Capita"
        $ws.Cells.Item($row, 4).Value = [double]$year
    }

    $ws.Cells.Item($row, 5).Value = $gdpPerCapita[$i]
}

Write-Output "Refreshed GDP per Capita data for rows $firstRow-$lastRow on the Data sheet"
